$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows at 391, pushing the existing 391-468 data block
# down to 394-471 (preserving all of its values/formatting).
$ws.Rows("391:393").Insert()

# Populate the 3 newly inserted rows with the new weekly records.

# Row 391
$ws.Range("A391").Value = 5
$ws.Range("B391").Value = "Macroferia Regional de Talca"
$ws.Range("C391").Value = "Maule"
$ws.Range("D391").Value = 44637
$ws.Range("E391").Value = 7
$ws.Range("F391").Value = 100112002
$ws.Range("G391").Value = "Pimiento"
$ws.Range("H391").Value = "Cuatro cascos rojo"
$ws.Range("I391").Value = "Primera"
$ws.Range("J391").Value = 200
$ws.Range("K391").Value = 15000
$ws.Range("L391").Value = 15000
$ws.Range("M391").Value = 15000
$ws.Range("N391").Value = '$/caja 15 kilos'
$ws.Range("O391").Value = "Región del Maule"
$ws.Range("P391").Value = 1000
$ws.Range("Q391").Value = 15
$ws.Range("R391").Value = "Hortaliza"

# Row 392
$ws.Range("A392").Value = 5
$ws.Range("B392").Value = "Macroferia Regional de Talca"
$ws.Range("C392").Value = "Maule"
$ws.Range("D392").Value = 44637
$ws.Range("E392").Value = 7
$ws.Range("F392").Value = 100112002
$ws.Range("G392").Value = "Pimiento"
$ws.Range("H392").Value = "Cuatro cascos verde"
$ws.Range("I392").Value = "Primera"
$ws.Range("J392").Value = 300
$ws.Range("K392").Value = 8000
$ws.Range("L392").Value = 8000
$ws.Range("M392").Value = 8000
$ws.Range("N392").Value = '$/caja 15 kilos'
$ws.Range("O392").Value = "Región del Maule"
$ws.Range("P392").Value = 533
$ws.Range("Q392").Value = 15
$ws.Range("R392").Value = "Hortaliza"

# Row 393
$ws.Range("A393").Value = 5
$ws.Range("B393").Value = "Macroferia Regional de Talca"
$ws.Range("C393").Value = "Maule"
$ws.Range("D393").Value = 44637
$ws.Range("E393").Value = 7
$ws.Range("F393").Value = 100112002
$ws.Range("G393").Value = "Pimiento"
$ws.Range("H393").Value = "Zafiro rojo"
$ws.Range("I393").Value = "Primera"
$ws.Range("J393").Value = 200
$ws.Range("K393").Value = 18000
$ws.Range("L393").Value = 18000
$ws.Range("M393").Value = 18000
$ws.Range("N393").Value = '$/caja 15 kilos'
$ws.Range("O393").Value = "Región de Arica y Parinacota"
$ws.Range("P393").Value = 1200
$ws.Range("Q393").Value = 15
$ws.Range("R393").Value = "Hortaliza"
